# "changes into switch window control" commit: refresh the AssetProfile
# (sheet2) asset-tag list.
#
#   rows 2-22   -> Auto-501 .. Auto-521   (existing rows, value only)
#   rows 23-41  -> asset tag cell removed entirely (kept their Cost cell)
#   rows 42-51  -> Auto-441 .. Auto-450   (existing rows, value only)
#   row  52     -> becomes a normal data row (was a distinct "total" style),
#                  Auto-451, Cost 1041
#   rows 53-59  -> brand-new data rows, Auto-452..Auto-458, Cost 1042-1048
#   row  60     -> brand-new row, Auto-459, no Cost
#
# Values are always written before formats are copied onto a cell, so a
# freshly created cell keeps its intended (numeric / string) type instead of
# inheriting the "Text" display format of the donor cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AssetProfile")

# --- rows 2-22: Auto-501 .. Auto-521 -------------------------------------
$n = 501
for ($r = 2; $r -le 22; $r++) {
    $ws.Cells.Item($r, 1).Value = "Auto-$n"
    $n++
}

# --- rows 23-41: remove the asset-tag cell entirely ----------------------
$ws.Range("A23:A41").Clear()

# --- rows 42-51: Auto-441 .. Auto-450 (existing rows, value only) -------
$n = 441
for ($r = 42; $r -le 51; $r++) {
    $ws.Cells.Item($r, 1).Value = "Auto-$n"
    $n++
}

# --- row 52 becomes a normal data row + rows 53-59 are brand-new data rows
# --- row 60 is a brand-new row with only an asset tag --------------------
$n = 451
$cost = 1041
for ($r = 52; $r -le 59; $r++) {
    $ws.Cells.Item($r, 1).Value = "Auto-$n"
    $ws.Cells.Item($r, 12).Value = $cost
    $n++
    $cost++
}
$ws.Cells.Item(60, 1).Value = "Auto-459"

# Now that every cell holds its final value, stamp them with the same
# look as the rest of the table (copy format only, so the stored value
# keeps its original type).
$ws.Range("A51").Copy()
$ws.Range("A52:A60").PasteSpecial(-4122)
$ws.Range("L51").Copy()
$ws.Range("L52:L59").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- selection moves to C14 ------------------------------------------------
$ws.Range("C14").Select()

Write-Output "AssetProfile asset list refreshed"
